$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Rules")

# Update cell B11 to the new text value "1" (it must remain a text/string
# cell, not get auto-coerced to a number, and it must keep its existing
# cell style). Writing a text-producing formula and then pasting back just
# the resulting value (Paste Values) achieves this without touching the
# cell's number format/style, unlike assigning .Value/.NumberFormat
# directly (which would re-point the cell at a freshly minted style).
$cell = $ws.Range("B11")
$cell.Formula = "=""1"""
$cell.Copy()
$cell.PasteSpecial(-4163)  # xlPasteValues
$excel.CutCopyMode = $false

$wb.Save()
